# Optuna Attempt (go back with original)
# Reverts forecast figures on "Forecast Comparison" and recomputed summary
# totals on "Summary" back to their original (pre-tweak) values.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison: MyForecast (D), Inventory Coverage (H),
#     Seasonality Index (L) ---
$ws1.Range("D2").Value  = 103
$ws1.Range("H2").Value  = 7.43
$ws1.Range("L2").Value  = 1.03

$ws1.Range("D3").Value  = 109
$ws1.Range("H3").Value  = 6.07
$ws1.Range("L3").Value  = 1.14

$ws1.Range("H4").Value  = 4.89
$ws1.Range("L4").Value  = 1.19

$ws1.Range("H5").Value  = 3.85
$ws1.Range("L5").Value  = 1.03

$ws1.Range("H6").Value  = 2.85
$ws1.Range("L6").Value  = 0.9399999999999999

$ws1.Range("H7").Value  = 1.89
$ws1.Range("L7").Value  = 0.84

$ws1.Range("H8").Value  = 0.88
$ws1.Range("L8").Value  = 1.02

$ws1.Range("L9").Value  = 0.96

$ws1.Range("D10").Value = 107
$ws1.Range("L10").Value = 1.03

$ws1.Range("L11").Value = 1.19

$ws1.Range("L12").Value = 1.2

$ws1.Range("L13").Value = 1.18

$ws1.Range("L14").Value = 1

$ws1.Range("D15").Value = 87
$ws1.Range("L15").Value = 1.09

$ws1.Range("D16").Value = 61
$ws1.Range("L16").Value = 1.11

$ws1.Range("D17").Value = 55
$ws1.Range("L17").Value = 0.88

# --- Summary: totals stored as text labels, keep them text (not numeric).
#     Force text entry via a "@" number format so Excel doesn't coerce the
#     digit string into a real number, then drop the format again so the
#     cell's style index is untouched (matches the original plain styling).
$ws2.Range("B9").NumberFormat  = "@"
$ws2.Range("B9").Value  = "1644"
$ws2.Range("B9").ClearFormats()

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "894"
$ws2.Range("B10").ClearFormats()

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "440"
$ws2.Range("B11").ClearFormats()

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "55"
$ws2.Range("B14").ClearFormats()
